# "modified for Selenium Grid"
#
# The AppFlow test-data sheet is trimmed down:
#   - the "Execute" result column (E) is dropped entirely
#   - only the first two data rows are kept (iterations 1 and 2)
#   - the browser for the second remaining row is switched to FIREFOX
#     (so the grid now exercises one Chrome node and one Firefox node)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Execute"/result column (column E).
$ws.Columns.Item(5).Delete() | Out-Null

# Drop the old iteration 3 & 4 rows (now rows 4 and 5).
$ws.Range("A4:D5").EntireRow.Delete() | Out-Null

# Point the remaining second data row (iteration 2) at Firefox instead of Chrome.
$ws.Range("D3").Value = "FIREFOX"

# Match the saved selection/active cell of the edited workbook.
$ws.Range("B4").Select() | Out-Null
